# Black_Hills_State_University_Organizations.xlsx — restructure "Organizations" sheet
#
# Summary of the target change:
#  - Swap columns A (was "Organization Name") and B (was "Categories") so the
#    category now leads: A="Category", B="Organization Name".
#  - Rename several headers (Org URL->Organization Link, Image URL->Logo Link,
#    Phone->Phone Number, Website dropped, LinkedIn->Linkedin Link,
#    Instagram->Instagram Link, Facebook->Facebook Link, Twitter->Twitter Link).
#  - Drop the old "Website" column (H) entirely (its data is discarded) and
#    shift LinkedIn/Instagram/Facebook/Twitter one column to the left (H..K).
#  - Append two brand-new, currently-empty columns: L="Youtube Link" and
#    M="Tiktok Link".
#  - Resize columns to match (A/B swap their widths; new L/M get 14/13).
#
# We do this by reading the full old grid into memory, then writing only the
# cells whose value actually changes (rather than using Columns.Insert/
# Delete, which leaves stray zero-span <col> entries behind in this engine,
# or blanket-rewriting every cell, which silently drops already-blank
# inlineStr cells whose write is a same-value no-op). The header's bold/
# border/alignment style (cellXfs index 1) is preserved by copying it onto
# the new header cell with a format-only paste. Brand-new blank cells
# (column M everywhere, and column L in rows 28-33 whose old content is
# discarded) are forced to materialize in the saved file by touching .Style
# right after clearing .Value - writing "" alone is a no-op there and the
# cell is simply omitted from the output.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 33
$oldCols = 12   # A..L in the original sheet
$newCols = 13   # A..M in the restructured sheet

# ---- 1. Snapshot every existing cell (A1:L33) ----------------------------
$oldData = @()
for ($r = 1; $r -le $lastRow; $r++) {
    $rowVals = @()
    for ($c = 1; $c -le $oldCols; $c++) {
        $v = $ws.Cells.Item($r, $c).Value2
        if ($v -eq $null) { $v = "" }
        $rowVals += ,$v
    }
    $oldData += ,$rowVals
}

# ---- 2. New header row (A..M) ---------------------------------------------
$headers = @(
    "Category",
    "Organization Name",
    "Organization Link",
    "Logo Link",
    "Description",
    "Email",
    "Phone Number",
    "Linkedin Link",
    "Instagram Link",
    "Facebook Link",
    "Twitter Link",
    "Youtube Link",
    "Tiktok Link"
)

# Give the brand-new header cell M1 the same style as the rest of row 1
# (bold, bordered, centered/top) *before* writing its text, by copying A1's
# format onto it.
$ws.Cells.Item(1, 1).Copy() | Out-Null
$ws.Range("M1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

for ($c = 1; $c -le $newCols; $c++) {
    $ws.Cells.Item(1, $c).Value = $headers[$c - 1]
}

# ---- 3. Data rows (2..33) --------------------------------------------------
# Only touch cells whose value actually changes:
#   - A/B always swap (Category <-> Organization Name).
#   - H..K only change on rows that actually carried social-link data
#     (28-33); elsewhere old H..L were already blank, so new H..K (= old
#     I..L) are still blank and need no write.
#   - L only changes (old Twitter text -> blank) on rows 28-33.
#   - M is a brand-new column on every row - always force it blank.
for ($r = 2; $r -le $lastRow; $r++) {
    $old = $oldData[$r - 1]

    $ws.Cells.Item($r, 1).Value = $old[1]   # A = old Categories
    $ws.Cells.Item($r, 2).Value = $old[0]   # B = old Organization Name

    if ($old[7] -ne "" -or $old[8] -ne "" -or $old[9] -ne "" -or $old[10] -ne "" -or $old[11] -ne "") {
        $ws.Cells.Item($r, 8).Value  = $old[8]    # H = old LinkedIn
        $ws.Cells.Item($r, 9).Value  = $old[9]    # I = old Instagram
        $ws.Cells.Item($r, 10).Value = $old[10]   # J = old Facebook
        $ws.Cells.Item($r, 11).Value = $old[11]   # K = old Twitter

        $ws.Cells.Item($r, 12).Value = ""
        $ws.Cells.Item($r, 12).Style = "Normal"   # force the now-blank L to persist
    }

    # M is a brand-new column on every row - force it to persist as blank.
    $ws.Cells.Item($r, 13).Value = ""
    $ws.Cells.Item($r, 13).Style = "Normal"
}

# ---- 4. Column widths -------------------------------------------------------
# OOXML stored widths (target): A=14 B=30 C=50 D=50 E=50 F=33 G=14 H=50 I=43
#                                J=42 K=31 L=14 M=13
# The ColumnWidth COM property reads ~0.83 below the stored OOXML width in
# this engine, so subtract that offset to land on the exact target width.
$targetWidths = @(14, 30, 50, 50, 50, 33, 14, 50, 43, 42, 31, 14, 13)
for ($c = 1; $c -le $newCols; $c++) {
    $ws.Columns.Item($c).ColumnWidth = $targetWidths[$c - 1] - 0.83
}
